# Update the generated output numbers (column F on each sheet) to reflect
# the latest scrape results, as published to gh-pages.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1158
$ws1.Range("F7").Value  = 819
$ws1.Range("F10").Value = 2065
$ws1.Range("F12").Value = 240
$ws1.Range("F13").Value = 100
$ws1.Range("F15").Value = 134
$ws1.Range("F16").Value = 2087
$ws1.Range("F17").Value = 558
$ws1.Range("F18").Value = 9661
$ws1.Range("F19").Value = 939
$ws1.Range("F21").Value = 92
$ws1.Range("F23").Value = 31
$ws1.Range("F24").Value = 246

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 533
$ws2.Range("F11").Value = 7
$ws2.Range("F13").Value = 40

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5637
$ws3.Range("F4").Value = 420

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5637
$ws4.Range("F5").Value  = 420
$ws4.Range("F6").Value  = 533
$ws4.Range("F7").Value  = 1158
$ws4.Range("F12").Value = 819
$ws4.Range("F16").Value = 2065
$ws4.Range("F18").Value = 240
$ws4.Range("F20").Value = 100
$ws4.Range("F24").Value = 134
$ws4.Range("F26").Value = 7
$ws4.Range("F27").Value = 2087
$ws4.Range("F28").Value = 558
$ws4.Range("F30").Value = 40
$ws4.Range("F31").Value = 939
$ws4.Range("F33").Value = 92
$ws4.Range("F36").Value = 31
$ws4.Range("F39").Value = 246
